$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text from "Total cost" to "Cost"
$ws.Range("B1").Value = "Cost"

# Adjust column widths. Column B already has the correct stored width and
# is left untouched. Columns A and C are new custom widths; the values
# below are chosen so the resulting stored OOXML width -- which this
# runtime derives via round(ColumnWidth*6)+5 pixels -- lands as close as
# possible to the target pixel widths (116px / 90px).
$ws.Columns.Item(1).ColumnWidth = 15
$ws.Columns.Item(3).ColumnWidth = 11.333333333333334

# Update the active selection
$null = $ws.Range("K8").Select()
